$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment author shown in the UI / stored in comments1.xml <author>.
$excel.UserName = "Author"

# --- New row 17: Wedos domain + hosting entry ---
$ws.Range("A17").Value = 42901
# Copy the date number format from A16 (mm-dd-yy / style "1") instead of
# re-creating a custom number format.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Wedos - objednání domény + hostingu, spuštění demo verze, "

# --- F1: extra cost note, right aligned ---
$ws.Range("F1").Value = "plus 483,-"
$ws.Range("F1").HorizontalAlignment = [Microsoft.Office.Interop.Excel.Constants]::xlRight

# Column F width
$ws.Columns("F:F").ColumnWidth = 10.75

# --- Comment on F1 ---
$comment = $ws.Range("F1").AddComment()
$comment.Text("Author:" + [char]10 + "Za zaplacení domény + hostingu") | Out-Null

# --- Selection moves to A17 ---
$ws.Range("A17").Select() | Out-Null
